$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$idx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Docente(s)*") {
        $idx = $i
        break
    }
}

# Insert a new paragraph right after it containing the professor's info,
# styled as a bulleted list item (ListBullet), matching the surrounding
# "Docente(s) Responsável(eis)" list entries used elsewhere in the doc.
$srcRange = $d.Paragraphs($idx).Range
$insertionPoint = $d.Range($srcRange.End, $srcRange.End)
$insertionPoint.InsertAfter("7455355 - Robson da Silva Rocha`r")

$newPara = $d.Paragraphs($idx + 1)
$newPara.Style = "ListBullet"
